# Generate Report for Archive
#
# 1. Update the status text shared by the "zh-cn"/"de-de" status cells from
#    "Ready for handoff" to "In Translation" (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" column(s) that held that text (Overview columns E & F,
#    and column C on the zh-cn / de-de sheets) to match the new, shorter content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status values -------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns -------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
